# Updated cryptos list on Sun Mar 17 14:43:43 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.319.70"
$ws.Range("E2").Value = "  -0.79%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.587.25"
$ws.Range("E3").Value = "  -1.77%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
$ws.Range("D5").Value = "577.52"
$ws.Range("E5").Value = "  -3.29%  "

# Row 6 - Solana
$ws.Range("D6").Value = "191.13"
$ws.Range("E6").Value = "  +0.49%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.581.97"
$ws.Range("E7").Value = "  -1.68%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.23%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.39%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  -3.04%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "0.151"
$ws.Range("E11").Value = "  -1.48%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "54.61"
$ws.Range("E12").Value = "  -4.50%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +0.40%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -2.43%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.160.22"
$ws.Range("E15").Value = "  -1.81%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.590.98"
$ws.Range("E16").Value = "  -1.76%  "

# Row 17 - TRON
$ws.Range("D17").Value = "0.125"
$ws.Range("E17").Value = "  -0.95%  "

# Row 18 - Uniswap
$ws.Range("D18").Value = "12.32"
$ws.Range("E18").Value = "  -0.76%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "67.213.76"
$ws.Range("E19").Value = "  -0.62%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "18.34"
$ws.Range("E20").Value = "  -2.65%  "

# Row 21 - Polygon
$ws.Range("E21").Value = "  -3.50%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "402.65"
$ws.Range("E22").Value = "  +0.37%  "

# Row 23 - RenderToken
$ws.Range("D23").Value = "13.11"
$ws.Range("E23").Value = "  +18.52%  "

# Row 24 - PancakeSwap
$ws.Range("E24").Value = "  -3.94%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "85.68"
$ws.Range("E25").Value = "  -2.03%  "

# Row 26 - ImmutableX
$ws.Range("D26").Value = "2.94"
$ws.Range("E26").Value = "  -0.22%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "12.55"
$ws.Range("E27").Value = "  +0.60%  "

# Row 28 - LEO
$ws.Range("E28").Value = "  +0.79%  "

# Row 29 - Toncoin
$ws.Range("D29").Value = "3.81"
$ws.Range("E29").Value = "  +4.29%  "

# Row 30 - NEARProtocol
$ws.Range("D30").Value = "8.09"
$ws.Range("E30").Value = "  +11.03%  "

# Row 31 - Filecoin
$ws.Range("D31").Value = "9.11"
$ws.Range("E31").Value = "  -1.69%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "31.21"
$ws.Range("E32").Value = "  -1.82%  "

# Row 33 - Bittensor
$ws.Range("D33").Value = "667.67"
$ws.Range("E33").Value = "  +9.97%  "

# Row 34 - Cosmos
$ws.Range("D34").Value = "12.19"
$ws.Range("E34").Value = "  -0.74%  "

# Row 35 - Hedera
$ws.Range("E35").Value = "  -0.45%  "

# Row 36 - OKB
$ws.Range("D36").Value = "64.09"
$ws.Range("E36").Value = "  -2.63%  "

# Row 37 - InjectiveProtocol
$ws.Range("D37").Value = "42.66"
$ws.Range("E37").Value = "  -3.79%  "

# Row 38 - TheGraph
$ws.Range("D38").Value = "0.422"
$ws.Range("E38").Value = "  +7.45%  "

# Row 39 - Dai
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.14%  "

# Row 40 - PEPE (contains a subscript-3 unicode char, build via interpolation)
$sub3 = [char]0x2083
$pepePrefix = "0.0"
$pepeSuffix = "0789"
$ws.Range("D40").Value = "$pepePrefix$sub3$pepeSuffix"
$ws.Range("E40").Value = "  +2.58%  "

# Row 41/42 - ThetaToken and Fetch.AI swap places (identity + values)
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "2.87"
$ws.Range("E41").Value = "  +13.71%  "

$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").Value = "3.13"
$ws.Range("E42").Value = "  +8.45%  "

# Row 43 - Kaspa
$ws.Range("E43").Value = "  -1.25%  "

# Row 44 - Maker
$ws.Range("D44").Value = "3.142.31"
$ws.Range("E44").Value = "  +13.24%  "

# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value = "  -0.14%  "

# Row 46 - VeChain
$ws.Range("D46").Value = "0.0418"
$ws.Range("E46").Value = "  -1.66%  "

# Row 47 - Stellar
$ws.Range("E47").Value = "  -3.03%  "

# Row 48 - Monero
$ws.Range("D48").Value = "143.48"
$ws.Range("E48").Value = "  +0.12%  "

# Row 49 - ApeXProtocol
$ws.Range("E49").Value = "  -2.06%  "

# Row 50 - THORChain
$ws.Range("D50").Value = "8.65"
$ws.Range("E50").Value = "  -0.87%  "

# Row 51 - WEMIXToken
$ws.Range("E51").Value = "  -2.91%  "
